{"js": "// The document originally has a title paragraph followed by five\n// \"Chapter N: ...\" paragraphs, each containing a heading, a body line and\n// an \"(Expand ...)\" note separated by manual line breaks.\n//\n// The edit rewrites those five chapter paragraphs into plain prose\n// paragraphs, then repeats that same 4-paragraph block (body style / fuel\n// type / regional preferences / strengths & weaknesses) five more times,\n// for a total of 24 new paragraphs after the title.\n\nconst bodyParaA =\n  \"Cars can be categorized by body style. Sedans are four-door vehicles offering comfort and practicality. Coupes are sportier, while hatchbacks provide flexible cargo space. SUVs dominate global markets for their size and versatility. Convertibles emphasize style and leisure, and trucks are valued for utility and towing capabilities.\";\nconst bodyParaB =\n  \"Fuel type also defines vehicles. Gasoline cars remain common due to infrastructure. Diesel engines provide torque and efficiency, especially for trucks. Hybrids combine combustion engines with electric motors to balance efficiency. Electric vehicles are gaining adoption, offering zero tailpipe emissions. Hydrogen fuel cell cars, though rare, represent another pathway.\";\nconst bodyParaC =\n  \"Regional preferences highlight consumer diversity. In the U.S., pickup trucks and SUVs dominate sales. In Europe, compact cars are popular due to narrow roads and fuel costs. Japan\\u2019s kei cars thrive in urban settings with small dimensions. Emerging markets embrace affordable economy vehicles.\";\nconst bodyParaD =\n  \"Each type has strengths and weaknesses. SUVs offer comfort and space but consume more fuel. Electric vehicles provide sustainability but require charging infrastructure. The car market continues to diversify to meet different consumer needs.\";\n\nconst cycle = [bodyParaA, bodyParaB, bodyParaC, bodyParaD];\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 0 is the \"Types of Cars\" title \u2014 leave it untouched.\n// Paragraphs 1-5 are the five \"Chapter N: ...\" paragraphs being replaced.\nconst chapterParas = body.paragraphs.items.slice(1, 6);\n\n// Rewrite the five existing chapter paragraphs in place with the first\n// cycle (A, B, C, D, A) \u2014 this removes their <w:br/> runs and\n// \"(Expand ...)\" notes, replacing each with a single plain-text run.\nfor (let i = 0; i < chapterParas.length; i++) {\n  chapterParas[i].insertText(cycle[i % cycle.length], \"Replace\");\n}\nawait context.sync();\n\n// Append five more full cycles (B,C,D,A / B,C,D,A / ... ) after the last\n// rewritten paragraph so the block (A,B,C,D) appears six times in total.\nlet lastPara = chapterParas[chapterParas.length - 1];\nconst totalNew = 24;\nfor (let i = chapterParas.length; i < totalNew; i++) {\n  lastPara = lastPara.insertParagraph(cycle[i % cycle.length], \"After\");\n}\nawait context.sync();\n", "ps1": "# The document originally has a title paragraph followed by five\n# \"Chapter N: ...\" paragraphs, each containing a heading, a body line and\n# an \"(Expand ...)\" note separated by manual line breaks.\n#\n# The edit rewrites those five chapter paragraphs into plain prose\n# paragraphs, then repeats that same 4-paragraph block (body style / fuel\n# type / regional preferences / strengths & weaknesses) five more times,\n# for a total of 24 new paragraphs after the title.\n\n$d = $word.ActiveDocument\n\n$bodyParaA = \"Cars can be categorized by body style. Sedans are four-door vehicles offering comfort and practicality. Coupes are sportier, while hatchbacks provide flexible cargo space. SUVs dominate global markets for their size and versatility. Convertibles emphasize style and leisure, and trucks are valued for utility and towing capabilities.\"\n$bodyParaB = \"Fuel type also defines vehicles. Gasoline cars remain common due to infrastructure. Diesel engines provide torque and efficiency, especially for trucks. Hybrids combine combustion engines with electric motors to balance efficiency. Electric vehicles are gaining adoption, offering zero tailpipe emissions. Hydrogen fuel cell cars, though rare, represent another pathway.\"\n$bodyParaC = \"Regional preferences highlight consumer diversity. In the U.S., pickup trucks and SUVs dominate sales. In Europe, compact cars are popular due to narrow roads and fuel costs. Japan\u2019s kei cars thrive in urban settings with small dimensions. Emerging markets embrace affordable economy vehicles.\"\n$bodyParaD = \"Each type has strengths and weaknesses. SUVs offer comfort and space but consume more fuel. Electric vehicles provide sustainability but require charging infrastructure. The car market continues to diversify to meet different consumer needs.\"\n\n$cycle = @($bodyParaA, $bodyParaB, $bodyParaC, $bodyParaD)\n\n# Paragraph 1 is the \"Types of Cars\" title \u2014 leave it untouched.\n# Paragraphs 2-6 are the five \"Chapter N: ...\" paragraphs being replaced\n# (Word COM paragraph indices are 1-based).\nfor ($i = 0; $i -lt 5; $i++) {\n    $p = $d.Paragraphs.Item($i + 2)\n    $p.Range.Text = $cycle[$i % 4]\n}\n\n# Append five more full cycles after the last rewritten paragraph so the\n# block (A,B,C,D) appears six times in total (24 new paragraphs overall).\n$lastPara = $d.Paragraphs.Item(6)\nfor ($i = 5; $i -lt 24; $i++) {\n    $lastPara.Range.InsertParagraphAfter()\n    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n    $lastPara.Range.Text = $cycle[$i % 4]\n}\n"}
